$d = $word.ActiveDocument

$d.Content.Find.Execute(" fMRI5", $true, $false, $false, $false, $false,
                         $true, 1, $false, " fBlo2", 2)

$d.Content.Find.Execute("Request Type: Scan Request", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Request Type: Lab Request", 2)

$d.Content.Find.Execute("Location: Lower Pike Hallway Exit Lobby", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Location: Multifaith Chapel", 2)

$d.Content.Find.Execute("Requested Employee: Leshin, Laurie", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Requested Employee: Franklin, Abraham", 2)
